$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04941833333333333
$ws.Range("H2").Value = 0.148255
$ws.Range("I2").Value = 0.005167549122999764
$ws.Range("J2").Value = 0.005167549122999764
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 3.809635760594444
$ws.Range("R2").Value = 34.28672184535
$ws.Range("S2").Value = 0.001242184074416926
$ws.Range("T2").Value = 0.001242184074416926

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04941833333333333
$ws.Range("H3").Value = 0.148255
$ws.Range("I3").Value = 0.005167549122999764
$ws.Range("J3").Value = 0.005167549122999764
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 5.019916144951112
$ws.Range("R3").Value = 45.17924530456
$ws.Range("S3").Value = 0.001636812619900881
$ws.Range("T3").Value = 0.001636812619900881

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04941833333333333
$ws.Range("H4").Value = 0.148255
$ws.Range("I4").Value = 0.005167549122999764
$ws.Range("J4").Value = 0.005167549122999764
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 7.018727217537778
$ws.Range("R4").Value = 63.16854495784001
$ws.Range("S4").Value = 0.002288552428681957
$ws.Range("T4").Value = 0.002288552428681957

# Row 5
$ws.Range("I5").Value = 0.806706161560336
$ws.Range("J5").Value = 0.806706161560336
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 594.722288694589
$ws.Range("R5").Value = 5352.5005982513
$ws.Range("S5").Value = 0.1939173721956901
$ws.Range("T5").Value = 0.1939173721956901

# Row 6
$ws.Range("I6").Value = 0.806706161560336
$ws.Range("J6").Value = 0.806706161560336
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.2555228396217449
$ws.Range("T6").Value = 0.2555228396217449

# Row 7
$ws.Range("I7").Value = 0.806706161560336
$ws.Range("J7").Value = 0.806706161560336
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.357265949742901
$ws.Range("T7").Value = 0.357265949742901

# Row 8
$ws.Range("I8").Value = 0.1881262893166642
$ws.Range("J8").Value = 0.1881262893166643
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 138.6910162302744
$ws.Range("R8").Value = 1248.21914607247
$ws.Range("S8").Value = 0.04522211110257539
$ws.Range("T8").Value = 0.0452221111025754

# Row 9
$ws.Range("I9").Value = 0.1881262893166642
$ws.Range("J9").Value = 0.1881262893166643
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.05958869033641394
$ws.Range("T9").Value = 0.05958869033641394

# Row 10
$ws.Range("I10").Value = 0.1881262893166642
$ws.Range("J10").Value = 0.1881262893166643
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.0833154878776749
$ws.Range("T10").Value = 0.08331548787767491
